$wb = $excel.ActiveWorkbook

# --- OFF sheet: update row 3 (R) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 238
$wsOff.Range("C3").Value = 163
$wsOff.Range("D3").Value = 54
$wsOff.Range("E3").Value = 28
$wsOff.Range("F3").Value = 4

# --- DEF sheet: update row 3 (R) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 222
$wsDef.Range("C3").Value = 147
$wsDef.Range("D3").Value = 63
$wsDef.Range("E3").Value = 38

$wb.Save()
